$d = $word.ActiveDocument

# 1. Remove the existing "_GoBack" bookmark (it currently sits near the top of
#    the document, left over from a prior edit session).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Delete the "Goldman Sachs Virtual Engineering Program" course bullet line
#    (the whole run of text), leaving the (now empty) paragraph in place.
$rng = $d.Content
$found = $rng.Find.Execute("Goldman Sachs Virtual Engineering Program", $true, $false, $false,
                            $false, $false, $true, 1, $false, "", 2)

# 3. Word re-creates the "_GoBack" bookmark at the location of the last edit
#    (the now-collapsed range where the text used to be).
if ($found) {
    $d.Bookmarks.Add("_GoBack", $rng)
}

# 4. The built-in "Normal Table" style picks up the QuickStyle (w:qFormat) flag
#    as part of this same save.
$tableStyle = $d.Styles("Normal Table")
$tableStyle.QuickStyle = $true
